$wb = $excel.ActiveWorkbook

# The shared status string "Ready for handoff" becomes "In Translation" for
# every language sheet (it is reused for the zh-cn and de-de status cells on
# the Overview sheet, and for the Status column on each language sheet).
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus

# Narrow the status-related columns (stored OOXML "character" width target is
# 13.4101845877511, which is ColumnWidth + 5/6 in this engine's width model).
$newColumnWidth = 13.4101845877511 - (5 / 6)

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
